$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# -------------------------------------------------------------------
# Two new "week" blocks are appended below the existing data (rows
# 1-60). Each block reuses the exact visual structure of the very
# first block (rows 1-10: title / column headers / 6 member rows /
# 2-row summary), so we clone the formatting from there and then only
# fill in the new text.
# -------------------------------------------------------------------

$ws.Range("A1:D10").Copy()
$ws.Range("A61:D70").PasteSpecial(-4122)

$ws.Range("A1:D10").Copy()
$ws.Range("A71:D80").PasteSpecial(-4122)

# ---- Block 1 : rows 61-70 -----------------------------------------
$ws.Range("A61").Value = "日期：2018.10.18 第七周周四"

$ws.Range("A62").Value = "组员"
$ws.Range("B62").Value = "计划内容"
$ws.Range("C62").Value = "完成情况"
$ws.Range("D62").Value = "备注"

$members = "练富珊","黄成志","黄皓燊","郑嘉蔚","陈碧容","辛伟达"
for ($i = 0; $i -lt 6; $i++) {
    $r = 63 + $i
    $ws.Range("A$r").Value = $members[$i]
    $ws.Range("B$r").Value = "课堂评审"
    $ws.Range("C$r").Value = "已完成"
}

$ws.Range("A69").Value = "总结：演示自己的作品，学习他人的优点。"

# ---- Block 2 : rows 71-80 -------------------------------------------
$ws.Range("A71").Value = "日期：2018.10.22 第八周周一"

$ws.Range("A72").Value = "组员"
$ws.Range("B72").Value = "计划内容"
$ws.Range("C72").Value = "完成情况"
$ws.Range("D72").Value = "备注"

for ($i = 0; $i -lt 6; $i++) {
    $r = 73 + $i
    $ws.Range("A$r").Value = $members[$i]
}

$ws.Range("A79").Value = "总结："

# ---- Merges matching the existing per-block layout -----------------
$ws.Range("A61:D61").Merge()
$ws.Range("A69:D70").Merge()
$ws.Range("A71:D71").Merge()
$ws.Range("A79:D80").Merge()

# ---- View state: scrolled down, I74 selected -----------------------
$ws.Range("I74").Select()
